$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text value looks like a plain number (single decimal point).
# Excel auto-converts such text to a numeric value on assignment, which would
# lose formatting (e.g. trailing zeros) and change the cell type. Force these
# specific cells to Text format before the assignment, then restore the default
# "Normal" style afterwards so no stray formatting is left behind.
$protectedCells = @(
    "D5", "D6", "D10", "D11", "D12", "D15", "D17", "D20", "D21", "D22", "D23", "D24", "D26", "D31", "D32", "D33", "D34", "D37", "D38", "D39", "D40", "D42", "D43", "D44", "D45", "D47", "D48", "D49", "D51"
)

foreach ($addr in $protectedCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the updated Coin / Link / Price / Volume(1h) values row by row.

# Row 2
$ws.Range("D2").Value = '44.187.47'
$ws.Range("E2").Value = '  +4.46%  '

# Row 3
$ws.Range("D3").Value = '2.221.44'
$ws.Range("E3").Value = '  +2.04%  '

# Row 4
$ws.Range("E4").Value = '  -0.01%  '

# Row 5
$ws.Range("D5").Value = '260.10'
$ws.Range("E5").Value = '  +2.71%  '

# Row 6
$ws.Range("D6").Value = '83.22'
$ws.Range("E6").Value = '  +13.14%  '

# Row 7
$ws.Range("E7").Value = '  +3.21%  '

# Row 8
$ws.Range("E8").Value = '  -0.05%  '

# Row 9
$ws.Range("E9").Value = '  +4.18%  '

# Row 10
$ws.Range("D10").Value = '44.27'
$ws.Range("E10").Value = '  +7.97%  '

# Row 11
$ws.Range("D11").Value = '0.0934'
$ws.Range("E11").Value = '  +2.47%  '

# Row 12
$ws.Range("D12").Value = '7.07'
$ws.Range("E12").Value = '  +4.19%  '

# Row 13
$ws.Range("E13").Value = '  +2.77%  '

# Row 14
$ws.Range("D14").Value = '2.555.29'
$ws.Range("E14").Value = '  +2.03%  '

# Row 15
$ws.Range("D15").Value = '14.60'
$ws.Range("E15").Value = '  +2.76%  '

# Row 16
$ws.Range("D16").Value = '2.205.04'
$ws.Range("E16").Value = '  +1.18%  '

# Row 17
$ws.Range("D17").Value = '0.783'
$ws.Range("E17").Value = '  +2.33%  '

# Row 18
$ws.Range("D18").Value = '44.067.84'
$ws.Range("E18").Value = '  +4.54%  '

# Row 19
$ws.Range("E19").Value = '  +1.78%  '

# Row 20
$ws.Range("D20").Value = '71.38'
$ws.Range("E20").Value = '  +0.97%  '

# Row 21
$ws.Range("D21").Value = '6.02'
$ws.Range("E21").Value = '  +2.63%  '

# Row 22
$ws.Range("D22").Value = '2.37'
$ws.Range("E22").Value = '  +9.58%  '

# Row 23
$ws.Range("D23").Value = '233.25'
$ws.Range("E23").Value = '  +2.87%  '

# Row 24
$ws.Range("D24").Value = '9.33'
$ws.Range("E24").Value = '  -1.88%  '

# Row 25
$ws.Range("E25").Value = '  +0.10%  '

# Row 26
$ws.Range("D26").Value = '10.80'
$ws.Range("E26").Value = '  +3.14%  '

# Row 27
$ws.Range("E27").Value = '  +12.37%  '

# Row 28
$ws.Range("E28").Value = '  +1.30%  '

# Row 29
$ws.Range("E29").Value = '  +2.51%  '

# Row 30
$ws.Range("E30").Value = '  +0.06%  '

# Row 31
$ws.Range("D31").Value = '172.97'
$ws.Range("E31").Value = '  +2.53%  '

# Row 32
$ws.Range("D32").Value = '20.67'
$ws.Range("E32").Value = '  +3.32%  '

# Row 33
$ws.Range("D33").Value = '0.0881'
$ws.Range("E33").Value = '  +9.52%  '

# Row 34
$ws.Range("D34").Value = '5.34'
$ws.Range("E34").Value = '  +4.30%  '

# Row 35
$ws.Range("E35").Value = '  +8.04%  '

# Row 36
$ws.Range("E36").Value = '  +2.30%  '

# Row 37
$ws.Range("D37").Value = '0.0363'
$ws.Range("E37").Value = '  +9.09%  '

# Row 38
$ws.Range("D38").Value = '4.50'
$ws.Range("E38").Value = '  +6.57%  '

# Row 39
$ws.Range("D39").Value = '13.53'
$ws.Range("E39").Value = '  +13.21%  '

# Row 40
$ws.Range("D40").Value = '2.96'
$ws.Range("E40").Value = '  +20.83%  '

# Row 41
$ws.Range("E41").Value = '  +3.31%  '

# Row 42
$ws.Range("B42").Value = 'THORChain'
$ws.Range("C42").Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range("D42").Value = '5.58'
$ws.Range("E42").Value = '  +9.05%  '

# Row 43
$ws.Range("B43").Value = 'MultiversX'
$ws.Range("C43").Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range("D43").Value = '63.47'
$ws.Range("E43").Value = '  +7.49%  '

# Row 44
$ws.Range("D44").Value = '0.202'
$ws.Range("E44").Value = '  +3.27%  '

# Row 45
$ws.Range("D45").Value = '103.17'
$ws.Range("E45").Value = '  +0.75%  '

# Row 46
$ws.Range("E46").Value = '  +2.27%  '

# Row 47
$ws.Range("D47").Value = '8.32'
$ws.Range("E47").Value = '  +0.44%  '

# Row 48
$ws.Range("B48").Value = 'Stacks'
$ws.Range("C48").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D48").Value = '1.57'
$ws.Range("E48").Value = '  +29.41%  '

# Row 49
$ws.Range("B49").Value = 'ARBITRUM'
$ws.Range("C49").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D49").Value = '1.12'
$ws.Range("E49").Value = '  +3.52%  '

# Row 50
$ws.Range("E50").Value = '  +3.86%  '

# Row 51
$ws.Range("D51").Value = '0.444'
$ws.Range("E51").Value = '  -5.32%  '

# Restore the default cell style on the protected cells (keeps them as text
# while dropping the temporary Text number-format override).
foreach ($addr in $protectedCells) {
    $ws.Range($addr).Style = "Normal"
}
